$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, shifting rows 22..78 down to 23..79.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new market record.
$ws.Cells.Item(22, 1).Value = 6
$ws.Cells.Item(22, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(22, 3).Value = "Metropolitana"
$ws.Cells.Item(22, 4).Value = 44690
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = 100114007
$ws.Cells.Item(22, 7).Value = "Jengibre"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 10000
$ws.Cells.Item(22, 12).Value = 11000
$ws.Cells.Item(22, 13).Value = 10425
$ws.Cells.Item(22, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(22, 15).Value = "Perú"
$ws.Cells.Item(22, 16).Value = 802
$ws.Cells.Item(22, 17).Value = 13
$ws.Cells.Item(22, 18).Value = "Hortaliza"
